$d = $word.ActiveDocument

# --- Paragraph 1: "Proxy Review Sheet" -----------------------------------
# Change paragraph style from Heading1 to Normal, keep centered alignment,
# and restyle the paragraph mark to bold 24pt while the run text becomes
# bold 26pt.
$p1 = $d.Paragraphs(1)

$p1.Style = "Normal"
$p1.Range.ParagraphFormat.Alignment = 1

# paragraph-mark run properties: bold, sz/szCs = 48 (24pt)
$p1.Range.Font.Bold = 1
$p1.Range.Font.BoldBi = 1
$p1.Range.Font.Size = 24
$p1.Range.Font.SizeBi = 24

# Now push the actual run ("Proxy Review Sheet" text) to bold, sz/szCs = 52
# (26pt) using InsertXML on the exact text range, since that is the only
# reliable way to give the run different complex-script (...Cs) sizing than
# the paragraph mark.
$runRange = $d.Range(0, 19)
$runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>Proxy Review Sheet</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$runRange.InsertXML($runXml)

# --- Paragraph 2: "Heading 2" ---------------------------------------------
# This whole paragraph (style Heading2) is removed from the document.
$p2 = $d.Paragraphs(2)
$p2.Range.Delete()
